# Applies the "2024-06-30 -> 2024-07-01" homework-sheet update:
#  - the date heading, and
#  - every division fact in the 5-column practice table.
# Cells are addressed by (row, col) rather than plain text Find/Replace
# because a couple of the old values ("42÷3=14, 0") repeat in more than
# one cell but must become different new values.

$d = $word.ActiveDocument

# --- Heading date -----------------------------------------------------
$d.Content.Find.Execute("2024-06-30 Sunday", $true, $false, $false, $false, `
    $false, $true, 1, $false, "2024-07-01 Monday", 2)

# --- Practice table -----------------------------------------------------
$t = $d.Tables.Item(1)

$cellEdits = @(
    @{Row=1;  Col=1; Old="80÷9=8, 8";   New="48÷7=6, 6"},
    @{Row=1;  Col=2; Old="32÷4=8, 0";   New="97÷8=12, 1"},
    @{Row=1;  Col=3; Old="86÷2=43, 0";  New="29÷8=3, 5"},
    @{Row=1;  Col=4; Old="85÷7=12, 1";  New="62÷7=8, 6"},
    @{Row=1;  Col=5; Old="65÷2=32, 1";  New="58÷5=11, 3"},

    @{Row=5;  Col=1; Old="30÷9=3, 3";   New="56÷8=7, 0"},
    @{Row=5;  Col=2; Old="87÷9=9, 6";   New="97÷3=32, 1"},
    @{Row=5;  Col=3; Old="94÷5=18, 4";  New="86÷9=9, 5"},
    @{Row=5;  Col=4; Old="70÷2=35, 0";  New="13÷7=1, 6"},
    @{Row=5;  Col=5; Old="14÷9=1, 5";   New="57÷8=7, 1"},

    @{Row=9;  Col=1; Old="17÷9=1, 8";   New="68÷5=13, 3"},
    @{Row=9;  Col=2; Old="78÷3=26, 0";  New="94÷8=11, 6"},
    @{Row=9;  Col=3; Old="44÷8=5, 4";   New="55÷3=18, 1"},
    @{Row=9;  Col=4; Old="67÷4=16, 3";  New="27÷3=9, 0"},
    @{Row=9;  Col=5; Old="79÷3=26, 1";  New="52÷6=8, 4"},

    @{Row=13; Col=1; Old="37÷8=4, 5";   New="79÷2=39, 1"},
    @{Row=13; Col=2; Old="90÷9=10, 0";  New="40÷3=13, 1"},
    @{Row=13; Col=3; Old="42÷3=14, 0";  New="57÷8=7, 1"},
    @{Row=13; Col=4; Old="70÷9=7, 7";   New="88÷6=14, 4"},
    @{Row=13; Col=5; Old="60÷4=15, 0";  New="98÷3=32, 2"},

    @{Row=17; Col=1; Old="69÷2=34, 1";  New="58÷7=8, 2"},
    @{Row=17; Col=2; Old="42÷3=14, 0";  New="33÷7=4, 5"},
    @{Row=17; Col=3; Old="35÷8=4, 3";   New="34÷9=3, 7"},
    @{Row=17; Col=4; Old="56÷2=28, 0";  New="62÷6=10, 2"},
    @{Row=17; Col=5; Old="83÷9=9, 2";   New="89÷8=11, 1"}
)

foreach ($edit in $cellEdits) {
    $cell = $t.Cell($edit.Row, $edit.Col)
    $rng = $cell.Range
    # Replace:=1 (wdReplaceOne) -- the cell's Range is used only to locate
    # the single occurrence to change; wdReplaceAll searches beyond the
    # supplied range and would clobber the OTHER cells that happen to
    # share the same old text (e.g. the two "42÷3=14, 0" cells below).
    $rng.Find.Execute($edit.Old, $true, $false, $false, $false, $false, `
        $true, 1, $false, $edit.New, 1)
}
